# Apply the changes described by the commit "comit for testing excel sheet"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new test data (username / password sample columns E & F) ---
$ws.Range("E3").Value = "user name"
$ws.Range("F3").Value = 123445

$ws.Range("E4").Value = "fjhjf"
$ws.Range("F4").Value = 13442

$ws.Range("E5").Value = "cnjzc"
$ws.Range("F5").Value = 23456

$ws.Range("E6").Value = "xfbj"
$ws.Range("F6").Value = 64664

$ws.Range("E7").Value = "vcgh"
$ws.Range("F7").Value = 456789

# Row height picked up a slightly taller auto-fit height for rows 6 & 7
$ws.Rows.Item(6).RowHeight = 14.9
$ws.Rows.Item(7).RowHeight = 14.9

# --- Selection / active cell moves to F7, and Sheet1 becomes the active tab ---
$ws.Range("F7").Select()

# --- Defined names that pointed at the (soon to be removed) Settings sheet
#     now resolve to broken references, same as Excel/Calc would show
#     after the sheet they pointed to is deleted ---
$wb.Names.Item("Action_Keywords").RefersTo = "=#ref!!`$d`$2:`$d`$11"
$wb.Names.Item("Home_Page").RefersTo = "=#ref!!`$b`$2:`$b`$11"
$wb.Names.Item("Login_Page").RefersTo = "=#ref!!`$c`$2:`$c`$11"
$wb.Names.Item("Page_Name").RefersTo = "=#ref!!`$a`$2:`$a`$11"

# --- Remove the now unused "Settings" and "Test Cases" sheets ---
$wb.Worksheets.Item("Settings").Delete()
$wb.Worksheets.Item("Test Cases").Delete()
